$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'296.77"
$ws.Range("E2").Value = "'3.82%"
$ws.Range("D3").Value = "'41.45"
$ws.Range("E3").Value = "'3.42%"
$ws.Range("D4").Value = "'5.037"
$ws.Range("E4").Value = "'-0.08%"
$ws.Range("D5").Value = "'0.07460"
$ws.Range("E5").Value = "'2.59%"
$ws.Range("D6").Value = "'4.354"
$ws.Range("E6").Value = "'1.00%"
$ws.Range("D7").Value = "'1.577"
$ws.Range("E7").Value = "'4.35%"
$ws.Range("D8").Value = "'0.9290"
$ws.Range("E8").Value = "'1.45%"
$ws.Range("D10").Value = "'0.1184"
$ws.Range("E10").Value = "'-1.42%"
$ws.Range("D11").Value = "'0.1832"
$ws.Range("E11").Value = "'7.49%"
$ws.Range("D12").Value = "'0.08910"
$ws.Range("E12").Value = "'3.27%"
$ws.Range("D13").Value = "'0.04178"
$ws.Range("E13").Value = "'0.53%"
$ws.Range("E14").Value = "'0.03%"
$ws.Range("D15").Value = "'0.001289"
$ws.Range("E15").Value = "'0.76%"
$ws.Range("D16").Value = "'0.005947"
$ws.Range("E16").Value = "'-0.62%"
$ws.Range("D17").Value = "'3.345"
$ws.Range("E17").Value = "'-1.69%"
$ws.Range("D18").Value = "'0.3295"
$ws.Range("E18").Value = "'0.35%"
$ws.Range("D19").Value = "'7.853"
$ws.Range("E19").Value = "'0.26%"
$ws.Range("E20").Value = "'4.77%"
$ws.Range("E21").Value = "'2.74%"
$ws.Range("D22").Value = "'0.04030"
$ws.Range("E22").Value = "'4.89%"
$ws.Range("D23").Value = "'0.001263"
$ws.Range("E23").Value = "'-0.63%"
$ws.Range("D24").Value = "'0.003864"
$ws.Range("E24").Value = "'1.99%"
$ws.Range("D25").Value = "'0.0001229"
$ws.Range("E25").Value = "'-4.23%"
$ws.Range("D26").Value = "'0.0003721"
$ws.Range("E26").Value = "'-0.27%"
$ws.Range("D38").Value = "'0.02394"
$ws.Range("E38").Value = "'4.29%"
$ws.Range("D39").Value = "'0.05204"
$ws.Range("E39").Value = "'5.71%"
$ws.Range("D40").Value = "'0.006874"
$ws.Range("E40").Value = "'1.13%"
$ws.Range("D41").Value = "'0.007775"
$ws.Range("E41").Value = "'1.35%"
$ws.Range("E42").Value = "'4.33%"
$ws.Range("E43").Value = "'-0.21%"
$ws.Range("D44").Value = "'0.007169"
$ws.Range("E44").Value = "'3.64%"
$ws.Range("D45").Value = "'0.3223"
$ws.Range("E45").Value = "'4.21%"
$ws.Range("D46").Value = "'0.00006224"
$ws.Range("E46").Value = "'-2.97%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.34%"
$ws.Range("D48").Value = "'0.04642"
$ws.Range("E48").Value = "'-81.57%"
$ws.Range("D49").Value = "'0.004199"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.34%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.34%"
